$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from an existing header cell (AB1) onto the
# three new header cells so they share the bold/bordered/centered style
# (style index 1) instead of Excel minting a brand-new style.
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Season record is the same for every player row on this sheet (2-48).
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 29).Value = 65
    $ws.Cells.Item($r, 30).Value = 97
    $ws.Cells.Item($r, 31).Value = 0
}
